# docs: Daily Scrum Sprint 3 Dia 2
# Appends the "Daily Scrum - Dia 2" entry (and trailing blank lines) after the
# last paragraph of the document, using Range.InsertXML so the inserted
# WordprocessingML (runs, proofErr markers, run-level formatting, emoji run,
# lastRenderedPageBreak, etc.) matches the source content exactly.
$d = $word.ActiveDocument

# Create a fresh empty paragraph at the very end of the document body; its
# Range is then used as the InsertXML target so the new XML *replaces* an
# empty paragraph instead of being merged into the prior (non-empty) one.
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$tail.Collapse(0)
$target = $d.Paragraphs.Last.Range

$xml = '<w:p/><w:p/><w:p/><w:p/><w:p><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>Daily</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Scrum - Día 2</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Facilitador: Jorge Samuel Solano Dorantes (Scrum </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Master</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Alán Osmar Peña Polo (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Owner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Developer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Frontend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>¿Qué hice ayer?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Implementé exitosamente el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>watchdog</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de conexión serial (T1.1)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">El </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>thread</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> verifica la conexión cada 5 segundos</w:t></w:r></w:p><w:p><w:r><w:t>La reconexión automática funciona perfectamente</w:t></w:r></w:p><w:p><w:r><w:t>Los logs están capturando todos los eventos de desconexión/reconexión</w:t></w:r></w:p><w:p><w:r><w:t>Solo se detectó 1 desconexión durante las pruebas de ayer, se reconectó en 3 segundos</w:t></w:r></w:p><w:p><w:r><w:t>Avancé en el componente de visualización de resultados</w:t></w:r></w:p><w:p><w:r><w:t>Diseñé el modal responsive</w:t></w:r></w:p><w:p><w:r><w:t>Preparé la estructura de datos para recibir telemetría</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>¿Qué haré hoy?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Completar el componente de visualización de resultados en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>frontend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T2.4)</w:t></w:r></w:p><w:p><w:r><w:t>Implementar animaciones CSS para victoria/derrota</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Agregar efecto de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>confetti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para victorias</w:t></w:r></w:p><w:p><w:r><w:t>Mostrar obstáculos esquivados, tiempo sobrevivido, resultado</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Apoyar en las pruebas de renderizado de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T2.5)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Validar que los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> se ven correctamente en el LCD físico</w:t></w:r></w:p><w:p><w:r><w:t>Documentar casos de prueba visuales</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Impedimentos detectados:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Ninguno. El </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>watchdog</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> está funcionando perfectamente y el componente de visualización va por buen camino.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Jorge Samuel Solano Dorantes (Scrum </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Master</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Developer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Embebido)</w:t></w:r></w:p><w:p><w:r><w:t>¿Qué hice a</w:t></w:r><w:r><w:t>y</w:t></w:r><w:r><w:t>er?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Implementé exitosamente la escritura en CGRAM para personaje </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T1.2) </w:t></w:r></w:p><w:p><w:r><w:t>El personaje se almacena correctamente en dirección 0x00</w:t></w:r></w:p><w:p><w:r><w:t>Los 8 bytes se escriben sin errores</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Implementé la escritura en CGRAM para obstáculo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T1.3) </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/></w:rPr><w:t>✅</w:t></w:r></w:p><w:p><w:r><w:t>El obstáculo se almacena en dirección 0x01</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Verificado que ambos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> coexisten sin sobrescribirse</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Creé la estructura de datos para configuración de nivel (T1.4) </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Struct</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>ocupa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> 18 </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>bytes:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> 8 (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>character</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>) + 8 (obstacle) + 1 (type) + 1 (value)</w:t></w:r></w:p><w:p><w:r><w:t>Implementé la confirmación de carga exitosa (T1.5)</w:t></w:r></w:p><w:p><w:r><w:t>El PIC envía: {"</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>":"</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>loaded</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>", "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>character</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>":"ok", "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>obstacle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>":"ok", "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>goal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>":"ok"}</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>¿Qué haré hoy?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Implementar inicialización del juego con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T2.1)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Cargar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de CGRAM al iniciar el juego</w:t></w:r></w:p><w:p><w:r><w:t>Posicionar personaje en la primera columna del LCD</w:t></w:r></w:p><w:p><w:r><w:t>Implementar sistema de movimiento del personaje con botones (T2.2)</w:t></w:r></w:p><w:p><w:r><w:t>Configurar entrada de botones (2 botones para 2 carriles)</w:t></w:r></w:p><w:p><w:r><w:t>Actualizar posición vertical del personaje</w:t></w:r></w:p><w:p><w:r><w:t>Renderizar movimiento en LCD</w:t></w:r></w:p><w:p><w:r><w:t>Comenzar implementación de generación y movimiento de obstáculos (T2.3)</w:t></w:r></w:p><w:p><w:r><w:t>Diseñar sistema de generación aleatoria</w:t></w:r></w:p><w:p><w:r><w:t>Crear buffer de obstáculos activos</w:t></w:r></w:p><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t>Impedimentos detectados:</w:t></w:r></w:p><w:p><w:r><w:t>IMPEDIMENTO POTENCIAL: El sistema de movimiento de obstáculos podría consumir mucha RAM. Necesito calcular cuidadosamente el tamaño del buffer de obstáculos para no exceder los límites de memoria. Planeo comenzar con un buffer de 6 obstáculos y expandir solo si hay memoria disponible.</w:t></w:r></w:p><w:p/><w:p/>'
$target.InsertXML($xml)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
